{"js": "// Replace the multiplication-problem text in each table cell with its\n// new value. Each old problem string occurs exactly once in the document,\n// so a simple search+replace per pair is unambiguous and leaves all run\n// formatting (font, size, paragraph alignment, etc.) untouched.\nconst replacements = [\n  [\"67\u00d761=\", \"91\u00d775=\"],\n  [\"65\u00d798=\", \"55\u00d778=\"],\n  [\"18\u00d713=\", \"52\u00d789=\"],\n  [\"53\u00d724=\", \"64\u00d776=\"],\n  [\"22\u00d791=\", \"28\u00d715=\"],\n  [\"94\u00d714=\", \"95\u00d733=\"],\n  [\"88\u00d787=\", \"82\u00d794=\"],\n  [\"36\u00d727=\", \"55\u00d797=\"],\n  [\"96\u00d757=\", \"15\u00d783=\"],\n  [\"84\u00d711=\", \"68\u00d760=\"],\n  [\"67\u00d722=\", \"25\u00d772=\"],\n  [\"21\u00d733=\", \"76\u00d762=\"],\n  [\"75\u00d749=\", \"34\u00d784=\"],\n  [\"62\u00d731=\", \"13\u00d778=\"],\n  [\"83\u00d753=\", \"80\u00d731=\"],\n  [\"60\u00d738=\", \"39\u00d799=\"],\n  [\"47\u00d728=\", \"88\u00d785=\"],\n  [\"56\u00d763=\", \"53\u00d760=\"],\n  [\"20\u00d752=\", \"46\u00d719=\"],\n  [\"19\u00d712=\", \"46\u00d789=\"],\n  [\"41\u00d796=\", \"87\u00d787=\"],\n  [\"64\u00d752=\", \"64\u00d791=\"],\n  [\"97\u00d733=\", \"87\u00d769=\"],\n  [\"62\u00d717=\", \"82\u00d775=\"],\n  [\"93\u00d746=\", \"37\u00d736=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for \"${oldText}\"`);\n  }\n\n  // Replace only the first occurrence (each value is unique in this document).\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Replace the multiplication-problem text in each table cell with its new\n# value. Each old problem string occurs exactly once in the document, so a\n# plain Find/Replace per pair is unambiguous and leaves all run formatting\n# (font, size, paragraph alignment, etc.) untouched.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"67\u00d761=\", \"91\u00d775=\"),\n    @(\"65\u00d798=\", \"55\u00d778=\"),\n    @(\"18\u00d713=\", \"52\u00d789=\"),\n    @(\"53\u00d724=\", \"64\u00d776=\"),\n    @(\"22\u00d791=\", \"28\u00d715=\"),\n    @(\"94\u00d714=\", \"95\u00d733=\"),\n    @(\"88\u00d787=\", \"82\u00d794=\"),\n    @(\"36\u00d727=\", \"55\u00d797=\"),\n    @(\"96\u00d757=\", \"15\u00d783=\"),\n    @(\"84\u00d711=\", \"68\u00d760=\"),\n    @(\"67\u00d722=\", \"25\u00d772=\"),\n    @(\"21\u00d733=\", \"76\u00d762=\"),\n    @(\"75\u00d749=\", \"34\u00d784=\"),\n    @(\"62\u00d731=\", \"13\u00d778=\"),\n    @(\"83\u00d753=\", \"80\u00d731=\"),\n    @(\"60\u00d738=\", \"39\u00d799=\"),\n    @(\"47\u00d728=\", \"88\u00d785=\"),\n    @(\"56\u00d763=\", \"53\u00d760=\"),\n    @(\"20\u00d752=\", \"46\u00d719=\"),\n    @(\"19\u00d712=\", \"46\u00d789=\"),\n    @(\"41\u00d796=\", \"87\u00d787=\"),\n    @(\"64\u00d752=\", \"64\u00d791=\"),\n    @(\"97\u00d733=\", \"87\u00d769=\"),\n    @(\"62\u00d717=\", \"82\u00d775=\"),\n    @(\"93\u00d746=\", \"37\u00d736=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
